$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from row 15 down to the new row 16 so the new row matches
# the existing style pattern (style 11 for column A, style 15 for B-F).
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row height used by the other data rows (19pt).
$ws.Rows.Item(16).RowHeight = 19

# Populate the new reference sequence row. Values are set in the same
# order that keeps the shared-string table ordering consistent with the
# expected output (sequenceID, virus_name, virus_full_name, virus_family,
# virus_genus).
$ws.Range("A16").Value = "MG599939"
$ws.Range("C16").Value = "WSBV"
$ws.Range("B16").Value = "Wuhan sharpbelly bornavirus"
$ws.Range("D16").Value = "Bornaviridae"
$ws.Range("E16").Value = "Cultervirus"

# Reflect the author's final on-screen selection around the newly added row.
$ws.Range("A13:F16").Select()
